# Rebuild "In-Class Exercise" (column D) links to use lowercase "exercise_"
# in the file-path portion of the URLs (was "Exercise_"), matching the
# commit "build with lowercase exercises". Also update the active
# selection to C15 (previously D15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column D whose RestructuredText links reference an
# "Exercise_*" html page that must become "exercise_*".
$rowsToFix = 2..28

foreach ($r in $rowsToFix) {
    $cell = $ws.Cells.Item($r, 4)  # column D
    $val = $cell.Value()
    if ($null -ne $val -and $val -is [string] -and $val.Contains("Exercise_")) {
        $cell.Value = $val.Replace("Exercise_", "exercise_")
    }
}

# Update the active cell selection shown in the file from D15 to C15.
$ws.Range("C15").Select()
